{"js": "// Update the division-problem table: the five data rows (0, 4, 8, 12, 16 \u2014\n// i.e. every 4th row, the others being blank \"answer\" rows) each get their\n// cell text swapped out for the new set of problems from the commit.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// rowIndex -> new text for the 5 cells in that row, left to right.\nconst rowUpdates = {\n  0: [\"65\u00f77=\", \"45\u00f74=\", \"46\u00f78=\", \"10\u00f78=\", \"16\u00f74=\"],\n  4: [\"41\u00f79=\", \"58\u00f78=\", \"43\u00f72=\", \"54\u00f78=\", \"13\u00f74=\"],\n  8: [\"53\u00f77=\", \"15\u00f77=\", \"91\u00f77=\", \"39\u00f72=\", \"91\u00f78=\"],\n  12: [\"75\u00f79=\", \"60\u00f73=\", \"60\u00f79=\", \"86\u00f75=\", \"54\u00f78=\"],\n  16: [\"14\u00f74=\", \"88\u00f79=\", \"90\u00f76=\", \"73\u00f75=\", \"15\u00f72=\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = rowUpdates[rowIndex];\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = newValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: the five data rows (1, 5, 9, 13, 17 in\n# 1-based COM indexing -- i.e. every 4th row, the others being blank \"answer\"\n# rows) each get their cell text swapped out for the new set of problems\n# from the commit.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"65\u00f77=\", \"45\u00f74=\", \"46\u00f78=\", \"10\u00f78=\", \"16\u00f74=\")\n    5  = @(\"41\u00f79=\", \"58\u00f78=\", \"43\u00f72=\", \"54\u00f78=\", \"13\u00f74=\")\n    9  = @(\"53\u00f77=\", \"15\u00f77=\", \"91\u00f77=\", \"39\u00f72=\", \"91\u00f78=\")\n    13 = @(\"75\u00f79=\", \"60\u00f73=\", \"60\u00f79=\", \"86\u00f75=\", \"54\u00f78=\")\n    17 = @(\"14\u00f74=\", \"88\u00f79=\", \"90\u00f76=\", \"73\u00f75=\", \"15\u00f72=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $row = $t.Rows.Item($rowIndex)\n    $newValues = $rowUpdates[$rowIndex]\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $row.Cells.Item($c).Range.Text = $newValues[$c - 1]\n    }\n}\n"}
